# Fill in the missing "count" figures (column G) on the exam-bill sheet
# so that the dependent rate formulas in column I (and the grand total
# in I32) recalculate to their proper amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G9").Value  = 60    # প্রশ্নপত্র প্রণয়ন        -> I9  = K9*G9
$ws.Range("G12").Value = 60    # উত্তরপত্র পরীক্ষণ         -> I12 = IF(...)
$ws.Range("G14").Value = 61    # ক্লাস টেষ্ট                -> I14 = K14*H14*G14
$ws.Range("G16").Value = 27    # ব্যবহারিক                 -> I16 = G16*K16/H16
$ws.Range("G18").Value = 118   # সেন্ট্রাল ভাইভা            -> I18 = G18*K18/H18
$ws.Range("G26").Value = 1     # ইনভিজিলেশন                -> I26 = K26*G26

# I32 (SUM(I9:I31)) recalculates automatically from the formulas above.
